# Insert a new weekly price record for "Femacal de La Calera / Coquimbo -
# Zanahoria" right above the current row 325. All of the rows that used to
# live at 325..400 shift down by one (to 326..401), and the freshly
# inserted row 325 receives the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 325:400 down to 326:401, creating a blank row 325.
$ws.Rows.Item(325).Insert()

# Populate the newly inserted row 325 with the new record.
$ws.Range("A325").Value = 3
$ws.Range("B325").Value = "Femacal de La Calera"
$ws.Range("C325").Value = "Coquimbo"
$ws.Range("D325").Value = 44785
$ws.Range("E325").Value = 5
$ws.Range("F325").Value = 100114013
$ws.Range("G325").Value = "Zanahoria"
$ws.Range("H325").Value = "Sin especificar"
$ws.Range("I325").Value = "Primera"
$ws.Range("J325").Value = 480
$ws.Range("K325").Value = 11000
$ws.Range("L325").Value = 11500
$ws.Range("M325").Value = 11240
$ws.Range("N325").Value = "$/saco 20 kilos"
$ws.Range("O325").Value = "Provincia de Quillota"
$ws.Range("P325").Value = 562
$ws.Range("Q325").Value = 20
$ws.Range("R325").Value = "Hortaliza"
